$p = $ppt.ActivePresentation

# --- Slide 11: "Harsher filtering also had slightly negative effects with OV" ---
# Bold the word "negative" in the note textbox's second paragraph.
$s11 = $p.Slides.Item(11)
$tb11 = $s11.Shapes.Item(3)
$para11 = $tb11.TextFrame.TextRange.Paragraphs(2)
$text11 = $para11.Text
$idx11 = $text11.IndexOf("negative")
$word11 = $para11.Characters($idx11 + 1, 8)
$word11.Font.Bold = 1

# --- Slide 12: title change + bold "negative" in the note textbox ---
$s12 = $p.Slides.Item(12)

# Title: "EFFECTs of Harsher Filtering" -> "EFFECTs of Different Smoothing"
$title12 = $s12.Shapes.Item(1)
$title12.TextFrame.TextRange.Text = "EFFECTs of Different Smoothing"

# Note: "0.1 smoothing with FV had slightly negative effects (compared to 0.01 smoothing)"
$tb12 = $s12.Shapes.Item(3)
$para12 = $tb12.TextFrame.TextRange.Paragraphs(2)
$text12 = $para12.Text
$idx12 = $text12.IndexOf("negative")
$word12 = $para12.Characters($idx12 + 1, 8)
$word12.Font.Bold = 1
